$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 144, shifting the existing rows 144:165 down to 145:166.
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with a new price-quote record
# (same product/quality/price series as the old row 144, new date 2021-11-22).
$ws.Range("A144").Value = 7
$ws.Range("B144").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C144").Value = "Ñuble"
$ws.Range("D144").Value = 44522
$ws.Range("E144").Value = 16
$ws.Range("F144").Value = "Fruta"
$ws.Range("G144").Value = 100104
$ws.Range("H144").Value = "Frutos de pepita"
$ws.Range("I144").Value = 100104005
$ws.Range("J144").Value = "Pera"
$ws.Range("K144").Value = "Packham's Triumph"
$ws.Range("L144").Value = "Primera"
$ws.Range("M144").Value = 120
$ws.Range("N144").Value = 10000
$ws.Range("O144").Value = 11000
$ws.Range("P144").Value = 10500
$ws.Range("Q144").Value = "$/caja 16 kilos empedrada"
$ws.Range("R144").Value = "Provincia de Curicó"
$ws.Range("S144").Value = 656
$ws.Range("T144").Value = 16
